# Rename header columns to use _FV2304 / _FV2310 suffixes instead of _old / _new,
# turn the data range into an Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the range into an Excel Table ("Table1")
$range = $ws.Range("A1:U60")
$listObject = $ws.ListObjects.Add(1, $range, $null, 1)
$listObject.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
